# Update column G ("K") values on Sheet1 rows 2-25 per the regenerated
# save_data (using K instead of Strike#, after recalculating std/mean and
# writing s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newK = @{
    2  = 1
    3  = 2
    4  = 2
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    10 = 0
    11 = 1
    12 = 0
    13 = 1
    14 = 1
    15 = 0
    16 = 1
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 0
    23 = 1
    24 = 1
    25 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
